$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 2105463.5
$ws.Range("I33").Value = 2105463.5
$ws.Range("K33").Value = 2105463.5
$ws.Range("M33").Value = -2105234.5
$ws.Range("H69").Value = 12750
$ws.Range("I69").Value = 8500
$ws.Range("J69").Value = 17000
$ws.Range("K69").Value = 25500
$ws.Range("L69").Value = 51000
$ws.Range("M69").Value = -24626
$ws.Range("N69").Value = -52748
$ws.Range("H72").Value = 12750
$ws.Range("I72").Value = 8500
$ws.Range("J72").Value = 17000
$ws.Range("K72").Value = 76500
$ws.Range("L72").Value = 153000
$ws.Range("M72").Value = -72132
$ws.Range("N72").Value = -161736
$ws.Range("H113").Value = 3640
$ws.Range("I113").Value = 3288
$ws.Range("J113").Value = 3710.4
$ws.Range("K113").Value = 3288
$ws.Range("L113").Value = 3710.4
$ws.Range("M113").Value = -34
$ws.Range("N113").Value = -10218.4
$ws.Range("H132").Value = 2793.85
$ws.Range("I132").Value = 1527.1333
$ws.Range("K132").Value = 4581.3999
$ws.Range("M132").Value = -2051.3999
$ws.Range("H137").Value = 4102.104
$ws.Range("J137").Value = 6942.1763
$ws.Range("L137").Value = 20826.5289
$ws.Range("N137").Value = -25926.5289
$ws.Range("H141").Value = 3515.2144
$ws.Range("I141").Value = 1690.8889
$ws.Range("J141").Value = 6799
$ws.Range("K141").Value = 5072.6667
$ws.Range("L141").Value = 20397
$ws.Range("M141").Value = 107.3333000000002
$ws.Range("N141").Value = -30757

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4075.3635
$ws.Range("I2").Value = 3118.5715
$ws.Range("K2").Value = 3118.5715
$ws.Range("M2").Value = -3005.5715
$ws.Range("H32").Value = 3370.8
$ws.Range("I32").Value = 2940.3215
$ws.Range("K32").Value = 2940.3215
$ws.Range("M32").Value = -2653.3215
$ws.Range("H61").Value = 6616
$ws.Range("I61").Value = 5732
$ws.Range("J61").Value = 7500
$ws.Range("K61").Value = 5732
$ws.Range("L61").Value = 7500
$ws.Range("M61").Value = -5520
$ws.Range("N61").Value = -7924
$ws.Range("H74").Value = 4449.154
$ws.Range("I74").Value = 2417.625
$ws.Range("J74").Value = 7699.6
$ws.Range("K74").Value = 2417.625
$ws.Range("L74").Value = 7699.6
$ws.Range("M74").Value = -1543.625
$ws.Range("N74").Value = -9447.6
$ws.Range("H77").Value = 4449.154
$ws.Range("I77").Value = 2417.625
$ws.Range("J77").Value = 7699.6
$ws.Range("K77").Value = 12088.125
$ws.Range("L77").Value = 38498
$ws.Range("M77").Value = -7720.125
$ws.Range("N77").Value = -47234
$ws.Range("H116").Value = 4075.3635
$ws.Range("I116").Value = 3118.5715
$ws.Range("K116").Value = 3118.5715
$ws.Range("M116").Value = -824.5715
$ws.Range("H122").Value = 4225.08
$ws.Range("I122").Value = 2884.6428
$ws.Range("J122").Value = 5931.091
$ws.Range("K122").Value = 8653.928400000001
$ws.Range("L122").Value = 17793.273
$ws.Range("M122").Value = -6203.928400000001
$ws.Range("N122").Value = -22693.273
$ws.Range("H132").Value = 7998.44
$ws.Range("I132").Value = 2920.4614
$ws.Range("K132").Value = 8761.3842
$ws.Range("M132").Value = -6231.3842
$ws.Range("H136").Value = 6616
$ws.Range("I136").Value = 5732
$ws.Range("J136").Value = 7500
$ws.Range("K136").Value = 17196
$ws.Range("L136").Value = 22500
$ws.Range("M136").Value = -14646
$ws.Range("N136").Value = -27600

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4075.3635
$ws.Range("I3").Value = 3118.5715
$ws.Range("K3").Value = 3118.5715
$ws.Range("M3").Value = -3004.5715
$ws.Range("H107").Value = 2154.8096
$ws.Range("I107").Value = 1262.6
$ws.Range("K107").Value = 1262.6
$ws.Range("M107").Value = 657.4000000000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 136669
$ws.Range("I86").Value = 200003.5
$ws.Range("K86").Value = 200003.5
$ws.Range("M86").Value = -198880.5
$ws.Range("H89").Value = 136669
$ws.Range("I89").Value = 200003.5
$ws.Range("K89").Value = 1000017.5
$ws.Range("M89").Value = -994401.5
$ws.Range("H122").Value = 4473.5454
$ws.Range("I122").Value = 3439
$ws.Range("K122").Value = 10317
$ws.Range("M122").Value = -7867
$ws.Range("H134").Value = 3920.0588
$ws.Range("I134").Value = 3393.2222
$ws.Range("K134").Value = 10179.6666
$ws.Range("M134").Value = -7644.6666

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 45877.5
$ws.Range("I5").Value = 89998.55499999999
$ws.Range("K5").Value = 269995.665
$ws.Range("M5").Value = -269883.665
$ws.Range("H33").Value = 2849024
$ws.Range("I33").Value = 2849024
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 17094144
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -17093861
$ws.Range("H50").Value = 550.4
$ws.Range("J50").Value = 1275
$ws.Range("L50").Value = 3825
$ws.Range("N50").Value = -4787
$ws.Range("H53").Value = 550.4
$ws.Range("J53").Value = 1275
$ws.Range("L53").Value = 3825
$ws.Range("N53").Value = -4787
$ws.Range("H56").Value = 6936
$ws.Range("I56").Value = 6936
$ws.Range("K56").Value = 6936
$ws.Range("M56").Value = -6406
$ws.Range("H80").Value = 53625
$ws.Range("J80").Value = 53625
$ws.Range("L80").Value = 160875
$ws.Range("N80").Value = -162747
$ws.Range("H83").Value = 53625
$ws.Range("J83").Value = 53625
$ws.Range("L83").Value = 482625
$ws.Range("N83").Value = -491985
$ws.Range("H107").Value = 62809.766
$ws.Range("J107").Value = 70777.07000000001
$ws.Range("L107").Value = 212331.21
$ws.Range("N107").Value = -216171.21
$ws.Range("H113").Value = 3367389.5
$ws.Range("I113").Value = 9259496
$ws.Range("J113").Value = 471.7143
$ws.Range("K113").Value = 27778488
$ws.Range("L113").Value = 1415.1429
$ws.Range("M113").Value = -27776318
$ws.Range("N113").Value = -5755.1429
$ws.Range("H114").Value = 1122.125
$ws.Range("I114").Value = 570.6667
$ws.Range("J114").Value = 2776.5
$ws.Range("K114").Value = 1712.0001
$ws.Range("L114").Value = 8329.5
$ws.Range("M114").Value = 1541.9999
$ws.Range("N114").Value = -14837.5
$ws.Range("H135").Value = 45877.5
$ws.Range("I135").Value = 89998.55499999999
$ws.Range("K135").Value = 809986.9949999999
$ws.Range("M135").Value = -807451.9949999999
$ws.Range("H136").Value = 7595.2
$ws.Range("I136").Value = 5993.3335
$ws.Range("J136").Value = 9998
$ws.Range("K136").Value = 17980.0005
$ws.Range("L136").Value = 29994
$ws.Range("M136").Value = -12880.0005
$ws.Range("N136").Value = -40194
$ws.Range("H138").Value = 1893
$ws.Range("I138").Value = 756.25
$ws.Range("J138").Value = 4166.5
$ws.Range("K138").Value = 2268.75
$ws.Range("L138").Value = 12499.5
$ws.Range("M138").Value = 2871.25
$ws.Range("N138").Value = -22779.5
$ws.Range("H139").Value = 1117.6666
$ws.Range("I139").Value = 1117.6666
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 3352.9998
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = 1787.0002
$ws.Range("N33").ClearContents()
$ws.Range("N139").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 9000
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 9000
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 9000
$ws.Range("N24").Value = -9346
$ws.Range("H46").Value = 17649.5
$ws.Range("J46").Value = 29299
$ws.Range("L46").Value = 29299
$ws.Range("N46").Value = -29611
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("H122").Value = 6119.778
$ws.Range("I122").Value = 6004.364
$ws.Range("J122").Value = 6301.143
$ws.Range("K122").Value = 18013.092
$ws.Range("L122").Value = 18903.429
$ws.Range("M122").Value = -15563.092
$ws.Range("N122").Value = -23803.429
$ws.Range("M24").ClearContents()
$ws.Range("N104").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 574.875
$ws.Range("J16").Value = 600
$ws.Range("L16").Value = 600
$ws.Range("N16").Value = -940
$ws.Range("H40").Value = 458457
$ws.Range("I40").Value = 669422.3
$ws.Range("K40").Value = 669422.3
$ws.Range("M40").Value = -669286.3
$ws.Range("H136").Value = 5247.1
$ws.Range("I136").Value = 3745.5
$ws.Range("J136").Value = 6248.1665
$ws.Range("K136").Value = 11236.5
$ws.Range("L136").Value = 18744.4995
$ws.Range("M136").Value = -8686.5
$ws.Range("N136").Value = -23844.4995

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 505000
$ws.Range("I7").Value = 10000
$ws.Range("K7").Value = 10000
$ws.Range("M7").Value = -9887
$ws.Range("H81").Value = 11646.615
$ws.Range("J81").Value = 16104.223
$ws.Range("L81").Value = 32208.446
$ws.Range("N81").Value = -34330.446
$ws.Range("H84").Value = 11646.615
$ws.Range("J84").Value = 16104.223
$ws.Range("L84").Value = 161042.23
$ws.Range("N84").Value = -171650.23
$ws.Range("H126").Value = 6356.125
$ws.Range("I126").Value = 7183.3335
$ws.Range("J126").Value = 5859.8
$ws.Range("K126").Value = 21550.0005
$ws.Range("L126").Value = 17579.4
$ws.Range("M126").Value = -19080.0005
$ws.Range("N126").Value = -22519.4
$ws.Range("H132").Value = 3926.0417
$ws.Range("I132").Value = 3096.394
$ws.Range("K132").Value = 9289.181999999999
$ws.Range("M132").Value = -6759.181999999999
$ws.Range("H136").Value = 305089.12
$ws.Range("I136").Value = 359109.5
$ws.Range("K136").Value = 1077328.5
$ws.Range("M136").Value = -1074778.5
